$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Refresh the "panel_query_time" timestamps on the "data" sheet (F2:F9) ---
$dataSheet.Range("F2").Value = "2021-10-05 14:22:39.672701"
$dataSheet.Range("F3").Value = "2021-10-05 14:22:39.672709"
$dataSheet.Range("F4").Value = "2021-10-05 14:22:39.672712"
$dataSheet.Range("F5").Value = "2021-10-05 14:22:39.672715"
$dataSheet.Range("F6").Value = "2021-10-05 14:22:39.672718"
$dataSheet.Range("F7").Value = "2021-10-05 14:22:39.672720"
$dataSheet.Range("F8").Value = "2021-10-05 14:22:39.672723"
$dataSheet.Range("F9").Value = "2021-10-05 14:22:39.672725"

# --- Add the new "metadata" worksheet, placed right after "data" ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$ws2.Name = "metadata"

# Re-use the bold + thin-border + center/top-aligned header style already defined
# in the workbook (style used by data!B1) instead of re-building fonts/borders by
# hand, so no duplicate style entries get minted.
$dataSheet.Range("B1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("A2").PasteSpecial(-4122)      # xlPasteFormats
$excel.CutCopyMode = $false

# Header row
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Data row
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Severe hypertriglyceridaemia"
$ws2.Range("C2").Value = 228

# D2 ("1.16") must stay a text value, not be coerced to a float. Build it as a
# text-formula result in a scratch cell and paste just the value in, which keeps
# the cell on the default (unstyled) format - unlike forcing a "@" text format or
# a leading apostrophe, both of which would mint a brand-new cell style.
$ws2.Range("Z1").Formula = "=""1.16"""
$ws2.Range("Z1").Copy()
$ws2.Range("D2").PasteSpecial(-4163)      # xlPasteValues
$ws2.Range("Z1").ClearContents()
$excel.CutCopyMode = $false

$ws2.Range("E2").Value = "2021-08-05T15:48:21.333628Z"
$ws2.Range("F2").Value = "2021-10-05 14:22:39.668935"
$ws2.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/228/?format=json"

# Keep "data" as the active sheet/tab, as it was before this edit.
$dataSheet.Activate()
[void]$dataSheet.Range("A1").Select()
